# "super dense all action implemented"
#
# Content edits only (this presentation also carries PowerPoint's own
# revision-tracking metadata in ppt/changesInfos/changesInfo1.xml and
# auto-recalculated `datetimeFigureOut` header/footer date fields in the
# slide layouts/master; those are environment/session side effects that
# aren't reachable - or meaningful - through the PowerPoint COM object
# model, so they're intentionally left alone here).

$p = $ppt.ActivePresentation

# ----------------------------------------------------------------------
# Slide 2 ("Evaluation Metrics:") -> title becomes "Evaluation:"
# ----------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$title2 = $s2.Shapes.Item(1)
$title2.TextFrame.TextRange.Text = "Evaluation:"

# Paragraph 3 of the content placeholder: "... Use 500 negative samples."
# becomes "... Use 700 negative samples." - retype just the "Use 500 "
# span (leaving the leading space and "negative samples" runs alone).
$body2 = $s2.Shapes.Item(2)
$para3 = $body2.TextFrame.TextRange.Paragraphs(3)
$useSpan = $para3.Characters(17, 8)
if ($useSpan.Text -eq "Use 500 ") {
    $useSpan.Text = "Use 700 "
}

# ----------------------------------------------------------------------
# Slide 6 ("Future:") -> update the first two bullets
# ----------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$body6 = $s6.Shapes.Item(2)
$tr6 = $body6.TextFrame.TextRange

$para1 = $tr6.Paragraphs(1)
$para1.Runs(1).Text = "Fix SAM (for later). "

$para2 = $tr6.Paragraphs(2)
$para2.Runs(1).Text = "Sampling ways in all action and dense all action. "
